# Regenerate s_val data to filter save games: update B2:G15 numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1.505614041169197,   1.65323645889881,      0.1529057820181812, 0.4998867070740569, 0, 3.811642989160245),
    @(0.1554434735375247,  0.000002220651329265522, 0.7127328510149897, 0.4998867070740569, 0, 1.368065252277901),
    @(0.06328177979961902, 0.3375848360084654,    0.1529057820181812, 0.4998867070740569, 1, 1.053659104900323),
    @(3.182878228561681,   1.65323645889881,      0.7127328510149897, 0.4998867070740569, 0, 6.048734245549538),
    @(0.1554434735375247,  0.3375848360084654,    3.082599426703578,  0.4998867070740569, 1, 4.075514443323626),
    @(3.182878228561681,   1.65323645889881,      0.1529057820181812, 0.4998867070740569, 0, 5.488907176552729),
    @(3.182878228561681,   1.65323645889881,      0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538),
    @(0.7287194209349384,  1.65323645889881,      0.1529057820181812, 0.4998867070740569, 1, 3.034748368925986),
    @(3.182878228561681,   1.65323645889881,      0.7127328510149897, 0.4998867070740569, 0, 6.048734245549538),
    @(0.1554434735375247,  0.3375848360084654,    0.7127328510149897, 0.4998867070740569, 1, 1.705647867635037),
    @(1.505614041169197,   1.65323645889881,      3.082599426703578,  0.4998867070740569, 1, 6.741336633845642),
    @(3.182878228561681,   1.65323645889881,      0.7127328510149897, 0.4998867070740569, 0, 6.048734245549538),
    @(3.182878228561681,   1.65323645889881,      0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538),
    @(1.505614041169197,   1.65323645889881,      0.7127328510149897, 0.4998867070740569, 1, 4.371470058157054)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $col = 2 + $j
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}
